# Generate Report for Handoff
# Adds a new tracked file (5d7baab8-2b3d-445c-a719-9c3245fc4841.md) as a new
# row to each of the three tables in the workbook: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$commit = "e20ea825643cdd36d6c5144e14f7801b2cb14bde"
$newFile = "5d7baab8-2b3d-445c-a719-9c3245fc4841.md"
$newFileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/$newFile"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("B3").Value = "e2e\" + $newFile
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newFileUrl, "", "", "e2e\" + $newFile) | Out-Null
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-02 06:46:50"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = $newFile
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $newFileUrl, "", "", $newFile) | Out-Null
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "5d7baab8-2b3d-445c-a719-9c3245fc4841.3e056a3a27c5c08329c7e26f3fba458a30ea3236.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-02 06:46:46"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("O3").Value = "False"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = $newFile
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $newFileUrl, "", "", $newFile) | Out-Null
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "5d7baab8-2b3d-445c-a719-9c3245fc4841.3e056a3a27c5c08329c7e26f3fba458a30ea3236.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-02 06:46:50"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("O3").Value = "False"

$wb.Save()
